# Update workbook to reflect data through 2021-10-21 (adds data for 2021-10-21)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2021-10-21"

# Update the label for October row
$ws.Range("A11").Value = "October (through 10-21)"

# Update October row (row 11) values for columns C:H (2016-2021)
$ws.Range("C11").Value = 33
$ws.Range("D11").Value = 42
$ws.Range("E11").Value = 50
$ws.Range("F11").Value = 34
$ws.Range("G11").Value = 101
$ws.Range("H11").Value = 129

# Update Total row (row 12) values for columns C:H (2016-2021)
$ws.Range("C12").Value = 462
$ws.Range("D12").Value = 669
$ws.Range("E12").Value = 598
$ws.Range("F12").Value = 456
$ws.Range("G12").Value = 1002
$ws.Range("H12").Value = 1376
